$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.528.61'
$ws.Range("E2").Value = '  +1.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.07'
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.77'
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6281'
$ws.Range("E6").Value = '  +1.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07427'
$ws.Range("E8").Value = '  +0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2953'
$ws.Range("E9").Value = '  +0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.45'
$ws.Range("E10").Value = '  +2.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07680'
$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.846.93'
$ws.Range("E12").Value = '  +0.57%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.033'
$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6792'
$ws.Range("E14").Value = '  +1.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.49'
$ws.Range("E15").Value = '  +1.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009163'
$ws.Range("E16").Value = '  +1.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.920'
$ws.Range("E17").Value = '  +0.74%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.527.25'
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.099.78'
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '245.15'
$ws.Range("E20").Value = '  +4.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.58'
$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.434'
$ws.Range("E23").Value = '  +3.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.55'
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1414'
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.561'
$ws.Range("E27").Value = '  +1.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.83'
$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06176'
$ws.Range("E29").Value = '  +11.06%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.499'
$ws.Range("E30").Value = '  +0.56%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.131'
$ws.Range("E31").Value = '  +0.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.105'
$ws.Range("E32").Value = '  +0.27%  '

$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.227'
$ws.Range("E33").Value = '  +1.52%  '

$ws.Range("E34").Value = '  +1.77%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7314'
$ws.Range("E35").Value = '  -1.57%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.146'
$ws.Range("E36").Value = '  +1.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.615'
$ws.Range("E37").Value = '  -1.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.893'
$ws.Range("E38").Value = '  +4.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.229.86'
$ws.Range("E39").Value = '  +2.03%  '

$ws.Range("E40").Value = '  -0.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.328'
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9169'
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("E43").Value = '  +0.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.014.40'
$ws.Range("E44").Value = '  +1.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.99'
$ws.Range("E45").Value = '  +0.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.88'
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000121'
$ws.Range("E47").Value = '  -0.61%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5067'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.285'
$ws.Range("E49").Value = '  +2.43%  '

$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4068'
$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1163'
$ws.Range("E51").Value = '  +5.70%  '
